$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: append a brand-new empty paragraph at the end of the document and
# return its Range (collapsed, ready for InsertXML to fill in content).
# ---------------------------------------------------------------------------
function New-TailParagraph {
    $count = $d.Paragraphs.Count
    $tail = $d.Paragraphs.Item($count).Range
    $tail.InsertParagraphAfter()
    return $d.Paragraphs.Item($d.Paragraphs.Count).Range
}

# Namespace declarations shared by every InsertXML fragment below.
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$drawNs = 'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing"'

function Add-TextParagraph {
    param([string]$InnerRunsXml)
    $target = New-TailParagraph
    $xml = '<w:p ' + $wNs + '>' + $InnerRunsXml + '</w:p>'
    $target.InsertXML($xml)
}

function Add-PictureParagraph {
    param(
        [string]$AnchorId,
        [string]$EditId,
        [string]$DocPrId,
        [string]$DocPrName,
        [string]$Cx,
        [string]$Cy,
        [string]$EffectExtent,
        [string]$RId,
        [bool]$LastRenderedPageBreak
    )
    $target = New-TailParagraph
    $lrpb = ""
    if ($LastRenderedPageBreak) { $lrpb = "<w:lastRenderedPageBreak/>" }
    $xml = '<w:p ' + $wNs + ' ' + $drawNs + '>' +
           '<w:r><w:rPr><w:noProof/></w:rPr>' + $lrpb +
           '<w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="' + $AnchorId + '" wp14:editId="' + $EditId + '">' +
           '<wp:extent cx="' + $Cx + '" cy="' + $Cy + '"/>' +
           '<wp:effectExtent ' + $EffectExtent + '/>' +
           '<wp:docPr id="' + $DocPrId + '" name="' + $DocPrName + '"/>' +
           '<wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr>' +
           '<a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
           '<pic:pic><pic:nvPicPr><pic:cNvPr id="1" name=""/><pic:cNvPicPr/></pic:nvPicPr>' +
           '<pic:blipFill><a:blip r:embed="' + $RId + '"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill>' +
           '<pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="' + $Cx + '" cy="' + $Cy + '"/></a:xfrm>' +
           '<a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic>' +
           '</a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>'
    $target.InsertXML($xml)
}

function Add-EmptyParagraph {
    $target = New-TailParagraph
    $target.InsertXML('<w:p ' + $wNs + '/>')
}

# ---------------------------------------------------------------------------
# 1) Remove everything from "Screen 1 -" (paragraph 8) through the end of the
#    body (paragraph 18, the last picture) — keeps paragraphs 1-7 which run
#    from "Working ScreenShots:" through "Part (B):" untouched.
# ---------------------------------------------------------------------------
$startPara = $d.Paragraphs.Item(8)
$endPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()

# ---------------------------------------------------------------------------
# 2) Rebuild the tail of the document in the new order described by the diff.
# ---------------------------------------------------------------------------

# "Brief Overview – "
Add-TextParagraph '<w:r><w:t xml:space="preserve">Brief Overview – </w:t></w:r>'

# Spring Boot overview paragraph
Add-TextParagraph '<w:r><w:t>The web application is built in Spring Boot, it has two major functionalities to identify duplicate records from CSV file and to classify using trained model weights as duplicates or non-duplicates.</w:t></w:r>'

# Perceptron / metaphone paragraph (with spell-check boundary markers around "metaphone")
Add-TextParagraph ('<w:r><w:t xml:space="preserve">I have built a single layer perceptron in this project, to calculate build the model, using weights calculated using </w:t></w:r>' +
                    '<w:proofErr w:type="spellStart"/><w:r><w:t>metaphone</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
                    '<w:r><w:t xml:space="preserve"> (phonetic similarity identifier library) output of the Company Names, and then transferred the model learning to identify, duplicates using emails and calculated the expected and observed weights which are similar but not yet accurately.</w:t></w:r>')

# Blank paragraph
Add-EmptyParagraph

# "JSON Object Accesses:" (single run, no bookmark now)
Add-TextParagraph '<w:r><w:t>JSON Object Accesses:</w:t></w:r>'

# Picture 4 (now referencing rId6), keeps its original anchor/edit ids
Add-PictureParagraph -AnchorId '4335E5CA' -EditId '2DF74614' -DocPrId '4' -DocPrName 'Picture 4' `
    -Cx '5943600' -Cy '3074035' -EffectExtent 'l="0" t="0" r="0" b="0"' -RId 'rId6' -LastRenderedPageBreak $true

# Two blank paragraphs
Add-EmptyParagraph
Add-EmptyParagraph

# Picture 5 (now referencing rId7), keeps its original anchor/edit ids
Add-PictureParagraph -AnchorId '17359C45' -EditId '07190B36' -DocPrId '5' -DocPrName 'Picture 5' `
    -Cx '5943600' -Cy '2411730' -EffectExtent 'l="0" t="0" r="0" b="7620"' -RId 'rId7' -LastRenderedPageBreak $false

# Blank paragraph
Add-EmptyParagraph

# "Screen 1 -"
Add-TextParagraph '<w:r><w:t>Screen 1 -</w:t></w:r>'

# Picture 2 (now referencing rId8, new anchor/edit ids)
Add-PictureParagraph -AnchorId '6F6A8241' -EditId '6413F2DE' -DocPrId '2' -DocPrName 'Picture 2' `
    -Cx '5943600' -Cy '3135630' -EffectExtent 'l="0" t="0" r="0" b="7620"' -RId 'rId8' -LastRenderedPageBreak $true

# Blank paragraph
Add-EmptyParagraph

# "Screen 2 - "
Add-TextParagraph '<w:r><w:t xml:space="preserve">Screen 2 - </w:t></w:r>'

# Picture 3 (now referencing rId9, new anchor/edit ids)
Add-PictureParagraph -AnchorId '6C592735' -EditId '3342000F' -DocPrId '3' -DocPrName 'Picture 3' `
    -Cx '5943600' -Cy '3151505' -EffectExtent 'l="0" t="0" r="0" b="0"' -RId 'rId9' -LastRenderedPageBreak $false

# Blank paragraph
Add-EmptyParagraph

# "To scale this up further, mapReduce could be used." + _GoBack bookmark
Add-TextParagraph ('<w:r><w:t xml:space="preserve">To scale this up further, </w:t></w:r>' +
                    '<w:proofErr w:type="spellStart"/><w:r><w:t>mapReduce</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
                    '<w:r><w:t xml:space="preserve"> could be used.</w:t></w:r>' +
                    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>')

# Two trailing blank paragraphs
Add-EmptyParagraph
Add-EmptyParagraph

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
